# "Add files via upload" - update the description text on the Metadata sheet.
#
# The only content-level change in the diff is the wording of the single
# shared string (cell A1 on the "Metadata" sheet): the sentence "Table
# describes the data used to make the figures" becomes "Table summarizes
# the data used to make the figures. The three PDF documents contains the
# raw data." The remaining diff hunks (fileVersion/rupBuild bump,
# xr:revisionPtr guid/coauth version, removal of the stale O19 selection)
# are Excel session/bookkeeping artifacts written automatically whenever the
# file is re-saved, so we just make sure the sheet is active with A1
# selected (its natural state after editing A1) rather than trying to poke
# those values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("A1").Value = "These data are used in Figures 4, 6, S2, and S4. They describe effect of timing of topical imidacloprid exposure on 5th and 6th instar corn earworm larvae. Table summarizes the data used to make the figures. The three PDF documents contains the raw data."

$ws.Activate()
$ws.Range("A1").Select()
